# Generate Report for handoff
# Updates the localization-status workbook:
#  - Overview / zh-cn / de-de "Status" cells change from
#    "Handoff transform failed" to "Ready for handoff"
#  - zh-cn / de-de sheets gain a "Latest Handoff File" hyperlink (column C)
#    pointing at the freshly produced .xlf handoff package
#  - "Latest Handoff Datetime" (column D) is stamped with the handoff time
#  - "Handoff Reason" (column H) switches from "Ignored" to "Include"

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob/02a010020c3e061c768bb6ef87bfdb16adc82630"
$mdName = "1d7d6837-711c-4007-a6ec-6372acb2f548.md"
$configName = ".localization-config"

function Set-LangSheet($sheetName, $xlfName, $handoffDatetime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status -> Ready for handoff
    $ws.Range("B2").Value = "Ready for handoff"

    # Handoff Reason -> Include
    $ws.Range("H2").Value = "Include"

    # Latest Handoff Datetime (plain text, matches existing column formatting)
    $ws.Range("D2").Value = $handoffDatetime

    # Rebuild the hyperlinks so that they end up in document order
    # (A2 source file, C2 new handoff file, A3 config file) - this mirrors
    # how Excel renumbers the r:id relationship ids on save.
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), "$repoBase/e2e/$mdName", [System.Type]::Missing, [System.Type]::Missing, $mdName) | Out-Null

    $ws.Range("C2").Value = $xlfName
    $ws.Hyperlinks.Add($ws.Range("C2"), "$repoBase/e2e/$xlfName", [System.Type]::Missing, [System.Type]::Missing, $xlfName) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A3"), "$repoBase/$configName", [System.Type]::Missing, [System.Type]::Missing, $configName) | Out-Null

    # Keep the original cornflower-blue underlined hyperlink look (the
    # workbook's pre-existing "HyperLink" cell style) on every linked cell.
    foreach ($addr in @("A2", "C2", "A3")) {
        $ws.Range($addr).Font.Underline = 2
        $ws.Range($addr).Font.Color = 15570276
    }
}

# Overview sheet: only the shared "Status" text changes (both language
# columns read off the same underlying status value).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"

Set-LangSheet "zh-cn" "1d7d6837-711c-4007-a6ec-6372acb2f548.e8bea255c400e53aa003e80092cd58c3ac3b34b0.zh-cn.xlf" "2016-01-08 20:18:34"
Set-LangSheet "de-de" "1d7d6837-711c-4007-a6ec-6372acb2f548.e8bea255c400e53aa003e80092cd58c3ac3b34b0.de-de.xlf" "2016-01-08 20:18:43"

Write-Output "Report generated for handoff"
